# Auto-generated Excel COM-interop script to apply cryptos list update
# Commit: Updated cryptos list on Wed Jan 10 16:42:57 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.445.92"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").Value = "2.427.18"
$ws.Range("E3").Value = "  +7.39%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.71"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.65"
$ws.Range("E6").Value = "  -5.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.82"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.804.74"
$ws.Range("E14").Value = "  +7.96%  "
$ws.Range("D15").Value = "2.429.62"
$ws.Range("E15").Value = "  +7.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.841"
$ws.Range("E16").Value = "  +6.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "45.382.03"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.39"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "0.0₃0941"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  +6.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.17"
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.06"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.79"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("E26").Value = "  +4.35%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.24"
$ws.Range("E28").Value = "  -9.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.25"
$ws.Range("E30").Value = "  +6.25%  "
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("E31").Value = "  +16.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.75"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.90"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0765"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +18.51%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.76"
$ws.Range("E39").Value = "  -8.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0297"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "2.018.86"
$ws.Range("E42").Value = "  +14.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.21"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.89"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.30"
$ws.Range("E46").Value = "  +26.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("E47").Value = "  -10.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  +9.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.02"
$ws.Range("E49").Value = "  +8.00%  "
$ws.Range("D50").Value = "2.672.51"
$ws.Range("E50").Value = "  +7.92%  "
$ws.Range("E51").Value = "  -0.95%  "
